$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Holidays 2019")

$ws.Range("C2").Value = '42 30 груз сер'
$ws.Range("G2").Value = '42, 30, груз, сер'
$ws.Range("C3").Value = 'б/к легк сер'
$ws.Range("G3").Value = '42, 30, груз, сер'
$ws.Range("C4").Value = 'б/к легк сер'
$ws.Range("G4").Value = 'б/к, легк, сер'
$ws.Range("G5").Value = 'б/к, легк, сер'
$ws.Range("C6").Value = 'легк сер'
$ws.Range("C7").Value = 'легк сер'
$ws.Range("G7").Value = 'легк, сер'
$ws.Range("C8").Value = 'легк сер'
$ws.Range("G8").Value = 'легк, сер'
$ws.Range("C9").Value = 'легк сер'
$ws.Range("G9").Value = 'легк, сер'
$ws.Range("C10").Value = '210B C H Type сер'
$ws.Range("G10").Value = 'легк, сер'
$ws.Range("C11").Value = 'груз LS-2 Type сер'
$ws.Range("G11").Value = '210B, C, H, Type, сер'
$ws.Range("C12").Value = '202B C Type сер'
$ws.Range("G12").Value = '210B, C, H, Type, сер'
$ws.Range("C13").Value = '202B C LS-2 H Type сер'
$ws.Range("G13").Value = 'груз, LS-2, Type, сер'
$ws.Range("C14").Value = 'груз б/к сер'
$ws.Range("G14").Value = '202B, C, Type, сер'
$ws.Range("C15").Value = 'б/к легк сер'
$ws.Range("G15").Value = '202B, C, LS-2, H, Type, сер'
$ws.Range("C16").Value = 'б/к легк сер'
$ws.Range("G16").Value = '202B, C, LS-2, H, Type, сер'
$ws.Range("C17").Value = '8 сх сер'
$ws.Range("G17").Value = '202B, C, LS-2, H, Type, сер'
$ws.Range("C18").Value = 'легк сер'
$ws.Range("G18").Value = 'груз, б/к, сер'
$ws.Range("C19").Value = 'легк сер'
$ws.Range("G19").Value = 'груз, б/к, сер'
$ws.Range("C20").Value = 'легк сер'
$ws.Range("G20").Value = 'груз, б/к, сер'
$ws.Range("C21").Value = 'легк сер'
$ws.Range("G21").Value = 'груз, б/к, сер'
$ws.Range("C22").Value = 'легк сер'
$ws.Range("G22").Value = 'б/к, легк, сер'
$ws.Range("G23").Value = 'б/к, легк, сер'
$ws.Range("G24").Value = 'легк, сер'
$ws.Range("G25").Value = 'легк, сер'
$ws.Range("G26").Value = 'легк, сер'
$ws.Range("G27").Value = 'легк, сер'
$ws.Range("G28").Value = 'легк, сер'
$ws.Range("G29").Value = 'легк, сер'
$ws.Range("G30").Value = 'легк, сер'
